$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add the two new column headers.
# Write order matters for shared-string table ordering, so write G1
# (heeftStopcontact) before F1 (masthoogte.standaardHoogte) to mirror the
# original authoring order.
$ws.Range("G1").Value = "heeftStopcontact"
$ws.Range("F1").Value = "masthoogte.standaardHoogte"

# The masthoogte.standaardHoogte values look like numbers ("10.00" etc.)
# but must be stored as text, matching the source data. Force the column
# to a text format before assigning the values, then restore the default
# "Normal" style so the cells are not left referencing a custom format.
$ws.Range("F2:F5").NumberFormat = "@"
$ws.Range("F2").Value = "10.00"
$ws.Range("F3").Value = "12.00"
$ws.Range("F4").Value = "18.00"
$ws.Range("F5").Value = "20.00"
$ws.Range("F2:F5").Style = "Normal"

# heeftStopcontact boolean flags for rows 4-6.
$ws.Range("G4").Value = $true
$ws.Range("G5").Value = $true
$ws.Range("G6").Value = $true

# Set the new column widths (character units) so the saved XML width is as
# close as possible to the authored best-fit widths (14 and ~27.86).
$ws.Columns.Item(5).ColumnWidth = 13.165
$ws.Columns.Item(6).ColumnWidth = 27

# Update the selected cell to match the final state of the workbook.
$ws.Range("F6").Select() | Out-Null
